$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Ravi"
$ws.Range("B2").Value = "Ravi@gmail.com"

$ws.Range("C5").Value = "Male 22"

$ws.Range("A5").Value = "RaviKUMAR 1"
$ws.Range("A6").Value = "RaviKUMAR 2"
$ws.Range("A7").Value = "RaviKUMAR 3"
$ws.Range("A8").Value = "RaviKUMAR 4"

$ws.Range("C6").Value = "Male 23"
$ws.Range("C7").Value = "Male 24"
$ws.Range("C8").Value = "Male 25"

$ws.Range("D6").Value = "A 3"
$ws.Range("D7").Value = "A 4"
$ws.Range("D8").Value = "A 5"

$ws.Range("C2").Value = "Ravi ppk"

$ws.Range("D5").Value = "A 2"

$ws.Range("B6").Value = 18
$ws.Range("B7").Value = 18
$ws.Range("B8").Value = 18
$ws.Range("C2").Select()
